$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2 through 13 from 2023-10-09 to 2023-10-13
$ws.Range("C2:C13").Value = "2023-10-13"
